$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.00214132762312634,0.194860813704497,0.993576017130621,0,0.993576017130621,0,0.91220556745182,0,0.991434689507495,0.931477516059957,0,0.00214132762312634,0.0406852248394004,0.00428265524625268,0.880085653104925,0.0107066381156317,0.00428265524625268,0.00214132762312634,0.00214132762312634,0.00214132762312634,0.00214132762312634,0.0107066381156317,0.00642398286937901),
    @(0.914346895074946,0.794432548179872,0.00428265524625268,0.933618843683084,0.00214132762312634,0.995717344753747,0.00214132762312634,0.997858672376874,0.00428265524625268,0.00214132762312634,0,0,0.93576017130621,0.00428265524625268,0.00642398286937901,0.967880085653105,0.019271948608137,0,0.989293361884368,0.989293361884368,0.991434689507495,0.0214132762312634,0.00214132762312634),
    @(0.00214132762312634,0.00428265524625268,0,0,0.00214132762312634,0.00214132762312634,0.0835117773019272,0.00214132762312634,0.00428265524625268,0.0663811563169165,1,0.997858672376874,0.00642398286937901,0.00428265524625268,0.00428265524625268,0,0.00428265524625268,0.00642398286937901,0.00428265524625268,0.00428265524625268,0.00428265524625268,0,0.985010706638116),
    @(0.0792291220556745,0.00642398286937901,0.00214132762312634,0.0663811563169165,0.00214132762312634,0.00214132762312634,0.00214132762312634,0,0,0,0,0,0.0171306209850107,0.987152034261242,0.107066381156317,0.0214132762312634,0.972162740899358,0.991434689507495,0.00214132762312634,0.00214132762312634,0.00214132762312634,0.967880085653105,0.00642398286937901)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $rowValues[$c]
    }
}
